$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Date, Publisher, Contact values on the Metadata sheet
$ws.Range("B8").Value = "2024-05-20T17:01:27+10:00"
$ws.Range("B9").Value = "D Foulkes - Northern Australia Regional Digital Health Collaborative"
$ws.Range("B10").Value = "D Foulkes - Northern Australia Regional Digital Health Collaborative (https://nardhc.org)"

# Insert a new row 11 for Jurisdiction / Australia, matching the formatting of the row below it
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "Australia"

Write-Host "done"
